# edit.ps1 - apply "remove cf flag; introduce cmp command instead" changes
#
# Summary of changes (see commit message / xml diff):
#  1. Update cached datetimeFigureOut field text (08.02.2018 -> 21.03.2018)
#     on the slide master and every slide layout.
#  2. Slide 2 ("ALU (integer) Instructions"):
#     a. Merge the title's two runs into one run (same text).
#     b. Bit-field numbering box: "10"-"5" -> "11 -"-"-"-" 5" (3 runs).
#     c. Small numbering textbox: "6" -> "7".
#     d. Merge the "R0 ... not used / cond. / flags" 3 runs into 1 run.
#     e. Remove the "11" bit-marker rectangle and its arrow.
#     f. Remove the "CF - change flags" description box.
#     g. Remove the extra "1" bit-marker textbox.
#     h. Reposition the "S - signed flag" description box and its arrow
#        to the space freed up by the removed shapes.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Update-DateField($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -like "*08.02.2018*") {
                $sh.TextFrame.TextRange.Text = "21.03.2018"
            }
        }
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. datetimeFigureOut cached text: 08.02.2018 -> 21.03.2018
#    (slide master + every slide layout)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
Update-DateField $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 2 edits
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)

# a. Title: merge "ALU (integer) " + "Instructions " -> one run.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "TEMP_TITLE_PLACEHOLDER"
$title.TextFrame.TextRange.Text = "ALU (integer) Instructions "

# b. "10"/"-"/"5" -> "11 -"/"-"/" 5"  (edit end-of-string run first so the
#    untouched offsets of earlier runs remain valid).
$bits105 = Get-ShapeById $s.Shapes 38
$tr = $bits105.TextFrame.TextRange
$tr.Characters(4, 1).Text = " 5"
$tr2 = $bits105.TextFrame.TextRange
$tr2.Characters(1, 2).Text = "11 -"

# c. "6" -> "7"
$six = Get-ShapeById $s.Shapes 39
$six.TextFrame.TextRange.Text = "7"

# d. Merge "R0 ... not used" + "cond. " + "flags" -> one run.
$r0row = Get-ShapeById $s.Shapes 42
$mergedText = "R0                 R1                 R2                                not used                cond. flags"
$r0row.TextFrame.TextRange.Text = "TEMP_R0_PLACEHOLDER"
$r0row.TextFrame.TextRange.Text = $mergedText

# e. Remove the "11" bit marker rectangle and its down-arrow.
(Get-ShapeById $s.Shapes 43).Delete()
(Get-ShapeById $s.Shapes 44).Delete()

# f. Remove the "CF - change flags" description textbox.
(Get-ShapeById $s.Shapes 47).Delete()

# g. Remove the second "1" bit marker textbox.
(Get-ShapeById $s.Shapes 49).Delete()

# h. Reposition the "S - signed flag" textbox (id 45) and its arrow
#    connector (id 46) into the space vacated above. EMU targets are
#    nudged by +0.5 EMU (in points) to counter the COM layer's
#    single-precision round-trip of Left/Top/Width/Height.
$emuPerPoint = 914400.0 / 72.0
function EmuToPt($emu) {
    return ($emu + 0.5) / $emuPerPoint
}

$sFlagBox = Get-ShapeById $s.Shapes 45
$sFlagBox.Left = EmuToPt 7208886
$sFlagBox.Top = EmuToPt 4271334

$sFlagArrow = Get-ShapeById $s.Shapes 46
$sFlagArrow.Left = EmuToPt 7340341
$sFlagArrow.Top = EmuToPt 2855235
